$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.152440309524536
$ws.Range("B1").Value = 4.206967353820801
$ws.Range("C1").Value = 2.196130514144897
$ws.Range("D1").Value = 1.69752824306488
$ws.Range("E1").Value = 1.536093831062317
